$wb = $excel.ActiveWorkbook

# Remove the "Currency" (col I) and "Expense Type" (col F) columns from both
# "Simple Fields" sheets, leaving A1:H2 (Receipt Name, Vendor Address,
# Phone Number, Receipt Date, Receipt Number, Tax Amount, Total Value, Items).
foreach ($name in @("Simple Fields", "Simple Fields - Formatted")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("I1:I2").EntireColumn.Delete()
    $ws.Range("F1:F2").EntireColumn.Delete()
}

# Update the garbled / obfuscated item description text in the "Items" table
# (column A, rows 2-12) on both "Items" sheets.
$itemDescriptionUpdates = @{
    "A2"  = "green onion Pancakes (1)"
    "A3"  = "Pan Fried Leek Dumplings #AT (2)"
    "A4"  = "Pork Xiao Long Bao(10) ¿*/ÅË#E(10)"
    "A5"  = "Q-BAO (5) #NEJ (5)"
    "A6"  = "Chicken potstickers *'ÈPJ`$3N5(6)"
    "A7"  = "Tomato Mushroom Steamed dumpli PEATTAMAMKE (6)"
    "A8"  = "Zucchini shrimp dumplings 7U#HA"
    "A9"  = "beef stew nodle soup (Non Spicy P¿#PJHE(T#)"
    "A10" = "dandan noodle #2H"
    "A11" = "banana naan bread EAA#"
    "A12" = "house made plum juice PUMgrt"
}

foreach ($name in @("Items", "Items - Formatted")) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($addr in $itemDescriptionUpdates.Keys) {
        $ws.Range($addr).Value = $itemDescriptionUpdates[$addr]
    }
}
